$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 11.80898833333333
$ws.Range("H2").Value = 35.426965
$ws.Range("I2").Value = 0.08059095716837197
$ws.Range("J2").Value = 0.08059095716837197
$ws.Range("M2").Value = 139.728498
$ws.Range("N2").Value = 419.185494
$ws.Range("O2").Value = 0.9065295391216045
$ws.Range("P2").Value = 0.9065295391216045
$ws.Range("Q2").Value = 1650.05220271619
$ws.Range("R2").Value = 14850.46982444571
$ws.Range("S2").Value = 0.07305808325921322
$ws.Range("T2").Value = 0.07305808325921322

$ws.Range("G3").Value = 11.80898833333333
$ws.Range("H3").Value = 35.426965
$ws.Range("I3").Value = 0.08059095716837197
$ws.Range("J3").Value = 0.08059095716837197
$ws.Range("O3").Value = 0.005362677585431591
$ws.Range("P3").Value = 0.005362677585431591
$ws.Range("Q3").Value = 9.761069640237224
$ws.Range("R3").Value = 87.84962676213502
$ws.Range("S3").Value = 0.0004321833195953058
$ws.Range("T3").Value = 0.0004321833195953058

$ws.Range("G4").Value = 11.80898833333333
$ws.Range("H4").Value = 35.426965
$ws.Range("I4").Value = 0.08059095716837197
$ws.Range("J4").Value = 0.08059095716837197
$ws.Range("O4").Value = 0.0881077832929639
$ws.Range("P4").Value = 0.0881077832929639
$ws.Range("Q4").Value = 160.3725368286028
$ws.Range("R4").Value = 1443.352831457425
$ws.Range("S4").Value = 0.007100690589563453
$ws.Range("T4").Value = 0.007100690589563453

$ws.Range("I5").Value = 0.8141849724511824
$ws.Range("J5").Value = 0.8141849724511824
$ws.Range("M5").Value = 139.728498
$ws.Range("N5").Value = 419.185494
$ws.Range("O5").Value = 0.9065295391216045
$ws.Range("P5").Value = 0.9065295391216045
$ws.Range("Q5").Value = 16669.95596546571
$ws.Range("R5").Value = 150029.6036891914
$ws.Range("S5").Value = 0.7380827278359066
$ws.Range("T5").Value = 0.7380827278359066

$ws.Range("I6").Value = 0.8141849724511824
$ws.Range("J6").Value = 0.8141849724511824
$ws.Range("O6").Value = 0.005362677585431591
$ws.Range("P6").Value = 0.005362677585431591
$ws.Range("Q6").Value = 98.61300194669425
$ws.Range("R6").Value = 887.5170175202481
$ws.Range("S6").Value = 0.004366211502159193
$ws.Range("T6").Value = 0.004366211502159193

$ws.Range("I7").Value = 0.8141849724511824
$ws.Range("J7").Value = 0.8141849724511824
$ws.Range("O7").Value = 0.0881077832929639
$ws.Range("P7").Value = 0.0881077832929639
$ws.Range("S7").Value = 0.07173603311311656
$ws.Range("T7").Value = 0.07173603311311656

$ws.Range("I8").Value = 0.1052240703804457
$ws.Range("J8").Value = 0.1052240703804457
$ws.Range("M8").Value = 139.728498
$ws.Range("N8").Value = 419.185494
$ws.Range("O8").Value = 0.9065295391216045
$ws.Range("P8").Value = 0.9065295391216045
$ws.Range("Q8").Value = 2154.400632657546
$ws.Range("R8").Value = 19389.60569391791
$ws.Range("S8").Value = 0.09538872802648468
$ws.Range("T8").Value = 0.09538872802648467

$ws.Range("I9").Value = 0.1052240703804457
$ws.Range("J9").Value = 0.1052240703804457
$ws.Range("O9").Value = 0.005362677585431591
$ws.Range("P9").Value = 0.005362677585431591
$ws.Range("S9").Value = 0.0005642827636770921
$ws.Range("T9").Value = 0.0005642827636770921

$ws.Range("I10").Value = 0.1052240703804457
$ws.Range("J10").Value = 0.1052240703804457
$ws.Range("O10").Value = 0.0881077832929639
$ws.Range("P10").Value = 0.0881077832929639
$ws.Range("S10").Value = 0.009271059590283888
$ws.Range("T10").Value = 0.009271059590283886

